{"js": "// Merge the split \"<id>p026r_1</id>\" runs into a single run, matching the\n// canonical edit: the paragraph had 3 runs (`<id>`, `p026r_1`, `</id>`) that\n// get combined into one run with the text `<id>p026r_1</id>`.\nconst body = context.document.body;\nconst paras = body.paragraphs;\nparas.load(\"items/text\");\nawait context.sync();\n\nlet target = null;\nfor (const p of paras.items) {\n  if (p.text === \"<id>p026r_1</id>\") {\n    target = p;\n    break;\n  }\n}\n\nif (target) {\n  // Replacing the whole paragraph range's text collapses the underlying runs\n  // into a single run (Word applies the formatting of the first original\n  // character to the replacement text), exactly mirroring the diff.\n  const range = target.getRange();\n  range.insertText(\"<id>p026r_1</id>\", \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Merge the split \"<id>p026r_1</id>\" runs into a single run, matching the\n# canonical edit: the paragraph had 3 runs (`<id>`, `p026r_1`, `</id>`) that\n# get combined into one run with the text `<id>p026r_1</id>`.\n$d = $word.ActiveDocument\n\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.Replacement.ClearFormatting()\n\n# FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n# MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace\n$rng.Find.Execute(\"<id>p026r_1</id>\", $false, $false, $false, $false, $false, $true, 1, $false, \"<id>p026r_1</id>\", 1) | Out-Null\n"}
